# Weekly update: insert a new Palta price record for "Tercera" quality
# (Hass, week of 2023-02-24) right before the existing row 108, shifting
# every subsequent record down by one row (dimension grows from T168 to
# T169), matching the "Fruta / hortaliza, semanal" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 108:168 down to 109:169 and open up a blank row 108.
$ws.Rows(108).Insert()

# Populate the newly inserted row with this week's record.
$ws.Cells.Item(108, 1).Value = 1
$ws.Cells.Item(108, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(108, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(108, 4).Value = 44981
$ws.Cells.Item(108, 5).Value = 15
$ws.Cells.Item(108, 6).Value = "Fruta"
$ws.Cells.Item(108, 7).Value = 100106
$ws.Cells.Item(108, 8).Value = "Oleaginosos"
$ws.Cells.Item(108, 9).Value = 100106002
$ws.Cells.Item(108, 10).Value = "Palta"
$ws.Cells.Item(108, 11).Value = "Hass"
$ws.Cells.Item(108, 12).Value = "Tercera"
$ws.Cells.Item(108, 13).Value = 400
$ws.Cells.Item(108, 14).Value = 25000
$ws.Cells.Item(108, 15).Value = 26000
$ws.Cells.Item(108, 16).Value = 25500
$ws.Cells.Item(108, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(108, 18).Value = "Perú"
$ws.Cells.Item(108, 19).Value = 2550
$ws.Cells.Item(108, 20).Value = 10
